# Apply updated crypto price/volume data (GitHub Actions symbol-list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '306.83'
    'E2' = '-0.71%'
    'D3' = '40.54'
    'E3' = '0.13%'
    'D4' = '5.050'
    'E4' = '-0.53%'
    'D5' = '0.07582'
    'E5' = '-2.64%'
    'D6' = '1.597'
    'E6' = '-4.05%'
    'E7' = '-4.32%'
    'D8' = '0.9031'
    'E8' = '-0.70%'
    'D9' = '0.1022'
    'E9' = '-1.62%'
    'D10' = '0.1758'
    'E10' = '0.15%'
    'D11' = '0.09064'
    'E11' = '0.72%'
    'D12' = '0.04240'
    'E12' = '-4.28%'
    'D13' = '0.1054'
    'E13' = '-0.26%'
    'D14' = '0.001240'
    'E14' = '-0.84%'
    'D15' = '0.005866'
    'E15' = '0.46%'
    'D16' = '3.351'
    'E16' = '-0.33%'
    'D17' = '4.269'
    'E17' = '-1.42%'
    'D18' = '0.3267'
    'E18' = '-2.97%'
    'D19' = '6.769'
    'E19' = '-5.09%'
    'D20' = '0.1360'
    'E20' = '-2.16%'
    'D21' = '0.2733'
    'E21' = '-4.19%'
    'D22' = '0.04186'
    'E22' = '0.34%'
    'D23' = '0.001229'
    'E23' = '0.76%'
    'D24' = '0.004058'
    'E24' = '-1.09%'
    'E25' = '6.44%'
    'D26' = '0.0003019'
    'E26' = '0.82%'
    'D38' = '0.02378'
    'E38' = '-1.21%'
    'D39' = '0.05158'
    'E39' = '-1.17%'
    'D40' = '0.007776'
    'E40' = '-2.46%'
    'D41' = '0.1287'
    'E41' = '-3.26%'
    'D42' = '0.007091'
    'E42' = '-6.47%'
    'D43' = '0.001925'
    'E43' = '-3.29%'
    'D44' = '0.008511'
    'E44' = '5.74%'
    'D45' = '0.3340'
    'E45' = '-0.61%'
    'D46' = '0.00006370'
    'E46' = '-5.50%'
    'D47' = '0.00000000753'
    'E47' = '-0.10%'
    'D48' = '0.004417'
    'E48' = '7.18%'
    'D49' = '0.006500'
    'E49' = '95.67%'
    'D50' = '0.00002109'
    'E50' = '-0.10%'
    'D51' = '0.0002008'
    'E51' = '-0.10%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}

